# Half finish code for the changes needed for new lipid classes addition
#
# Adds two new fatty-acid / lipid classes (FA25:1 and FA25:2) to the
# whitelist table on Sheet1, inserted right after the existing "FA25:0"
# row (row 54) and before "FA26:0" (old row 54, now pushed down).
# Both new rows get the same "T" marker as their neighbours in the
# fa1 / fa2 / fa3 / TG columns (A/B/C/D/F), matching the existing pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two fresh blank rows at position 54, pushing "FA26:0" and
# everything below it down by two rows.
$ws.Rows.Item(54).Insert()
$ws.Rows.Item(54).Insert()

# New row 54: FA25:1
$ws.Range("A54").Value = "FA25:1"
$ws.Range("B54").Value = "T"
$ws.Range("C54").Value = "T"
$ws.Range("D54").Value = "T"
$ws.Range("F54").Value = "T"

# New row 55: FA25:2
$ws.Range("A55").Value = "FA25:2"
$ws.Range("B55").Value = "T"
$ws.Range("C55").Value = "T"
$ws.Range("D55").Value = "T"
$ws.Range("F55").Value = "T"

# Leave the selection where the editor last left off.
[void]$ws.Range("F54").Select()
